# edit.ps1 -- apply the "added to the user guide...docx" commit.
#
# Summary of the change:
#   * The bullet "Project File (available on github 4shanob)" becomes
#     "Project File (available on github <hyperlink>)" where the
#     hyperlink text/target is https://github.com/4Shaneob/NewUITest.
#   * Word's "last edit" marker (the hidden _GoBack bookmark) moves from
#     its old position (inside "Continue through the mock payment until
#     complete.") to the end of the paragraph we just edited. Removing it
#     from the old spot also re-merges the two runs that its start/end
#     tags had split in two, back into a single contiguous run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Drop the stale _GoBack bookmark sitting inside the "mock payment"
#    bullet, then re-join the run that its start/end tags had split in
#    two, restoring a single contiguous run of text.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$mockPay = $d.Content
$mockPay.Find.Execute(
    "Continue through the mock payment until complete.", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "Continue through the mock payment until complete.", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Turn the placeholder "4shanob" into the real hyperlink to the repo.
# ---------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("4shanob", $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0) | Out-Null
$target.Text = "https://github.com/4Shaneob/NewUITest"
$d.Hyperlinks.Add($target, "https://github.com/4Shaneob/NewUITest") | Out-Null

# ---------------------------------------------------------------------
# 3) Re-plant _GoBack right after the closing ")" that follows the new
#    hyperlink -- i.e. at the point we just finished editing, matching
#    how Word itself tracks the last edit location. A temporary marker
#    word is typed in, bookmarked, and deleted again so the bookmark
#    lands as a clean, collapsed bookmark squarely between the ")" run
#    and the following paragraph mark (rather than swallowing the ")"
#    itself, which a directly-collapsed range fails to do reliably).
# ---------------------------------------------------------------------
$marker = "GoBackMarker"

$afterEdit = $d.Range($target.End, $d.Content.End)
$afterEdit.Find.Execute(")", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0) | Out-Null
$afterEdit.InsertAfter($marker)

$markerRng = $d.Content
$markerRng.Find.Execute($marker, $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $markerRng) | Out-Null

$clearRng = $d.Content
$clearRng.Find.Execute($marker, $true, $false, $false, $false, $false,
                        $true, 1, $false, "", 2) | Out-Null
